# Apply the "commands" worksheet update:
#  - insert two new rows (addlist, addentries) after the "options" row
#  - fix a couple of typos in existing rows
#  - append a new "makebase" row with centered / wrapped formatting
#  - update the active selection to the new last row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two rows for the new "addlist" / "addentries" commands -----
# (old rows 4,5,6 shift down to become rows 6,7,8)
$ws.Range("A4:A5").EntireRow.Insert()

$ws.Range("A4").Value = "addlist"
$ws.Range("B4").Value = "Dodaje liste"
$ws.Range("C4").Value = "nazwa listy"

$ws.Range("A5").Value = "addentries"
$ws.Range("B5").Value = "Dodaje wpisy do listy"
$ws.Range("C5").Value = "nazwa listy, ilosc powtórzeń wpisu, wpis"

# --- Fix typos in existing rows -----------------------------------------
$ws.Range("B3").Value = "Wypisuje zawartość listy"
$ws.Range("C3").Value = "nazwa listy, domyślnie main"

$ws.Range("B8").Value = "Pokazuje ilość punktów dla użytkownika"

# --- Add the new "makebase" row -----------------------------------------
$ws.Range("A9").Value = "makebase"
$ws.Range("B9").Value = "Dodaje folder użytkownikowi, jeśli folder ten wcześniej nie istniał"
$ws.Range("C9").Value = "brak"

$a9 = $ws.Range("A9")
$a9.HorizontalAlignment = -4108
$a9.VerticalAlignment = -4108
$a9.Copy()
$ws.Range("C9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$b9 = $ws.Range("B9")
$b9.WrapText = $true
$b9.HorizontalAlignment = -4108
$b9.VerticalAlignment = -4108

$ws.Rows(9).RowHeight = 36

# --- Update selection ----------------------------------------------------
$ws.Range("A9:C9").Select() | Out-Null
